$d = $word.ActiveDocument

# 1) The bookmark ("_GoBack") previously sat at the end of the "Changes"
#    paragraph. The new content appended below moves the last-edit point to
#    the very end of the document, so Word's automatic _GoBack bookmark
#    needs to move there too. Start by removing it from its old location.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 2) The document currently ends with two empty paragraphs. Drop the first
#    of the two, leaving a single trailing empty paragraph that will act as
#    the insertion point for all of the new log entries.
$paras = $d.Paragraphs
$trailingEmpty = $paras.Item(13)
$trailingEmpty.Range.Delete()

# 3) Insert all of the new log-book content (new dated entries, plus the
#    moved _GoBack bookmark at the very end) at the end of the document.
#    Using a fully collapsed range at the end of the last (empty) paragraph
#    lets InsertXML fold the first new paragraph into that empty paragraph
#    rather than leaving a stray blank line behind.
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
$insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newContentXml = '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>1</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>/07/2021</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">All Networking functions have been moved to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NetworkUtility</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Instead of functions depending on the cryptography, instead there will be functions for synchronous/asynchronous sending/receiving. It will be the responsibility of the client/server to build the transmission formats/interpret what they receive. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NetworkUtility</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> will only be concerned with how something is sent, not the content that is sent.</w:t></w:r><w:r><w:br/><w:t>An asynchronous socket connect method has been implemented, this way when trying to send a message to an offline contact the user is not blocked while waiting for the connection to fail.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Yet another class may be implemented, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TransmissionFormatter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. This class could be responsible for taking information and returning a formatted byte array ready for transmission, the array can then be passed to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NetworkUtility</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for either synchronous or asynchronous transmission.</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">It can also be used for decoding a transmission (e.g. with RSA it can extract and return a tuple </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>containing a byte array with the data and a byte array with the signature, which the server will then worry about validating).</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t xml:space="preserve">Implemented Asynchronous socket connection in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>NetworkUtility</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Later on</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">A new concern has been realised, when an exception occurs in a thread (e.g. sending/receiving data) the exception is </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>not</w:t></w:r><w:r><w:t xml:space="preserve"> thrown to the invoker if it is in another thread. As a result one of two things must happen, either each method is provided with a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>callback</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to use for exceptions, or I need to find something that allows that to happen, from superficial research I’ve seen that Tasks might allow that.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="ED7D31" w:themeColor="accent2"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="ED7D31" w:themeColor="accent2"/></w:rPr><w:t>Research Tasks</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>1</w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t>/07/2021</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">I think I understand enough about tasks to try to migrate every </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>async</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> function to a task. A test using tasks for the initial connection was successful.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t xml:space="preserve">New GIT branch: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>TaskMigration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$insertionPoint.InsertXML($newContentXml)
